$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3747116666666667
$ws.Range("H2").Value = 1.124135
$ws.Range("I2").Value = 0.3914669751594584
$ws.Range("J2").Value = 0.3914669751594584
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.030956000000001
$ws.Range("N2").Value = 18.092868
$ws.Range("O2").Value = 0.364814105361131
$ws.Range("P2").Value = 0.3648141053611309
$ws.Range("Q2").Value = 2.259869574353334
$ws.Range("R2").Value = 20.33882616918001
$ws.Range("S2").Value = 0.1428126743212259
$ws.Range("T2").Value = 0.1428126743212259

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3747116666666667
$ws.Range("H3").Value = 1.124135
$ws.Range("I3").Value = 0.3914669751594584
$ws.Range("J3").Value = 0.3914669751594584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.789877333333333
$ws.Range("N3").Value = 20.369632
$ws.Range("O3").Value = 0.4107214552505144
$ws.Range("P3").Value = 0.4107214552505143
$ws.Range("Q3").Value = 2.544246252035556
$ws.Range("R3").Value = 22.89821626832
$ws.Range("S3").Value = 0.1607838857200097
$ws.Range("T3").Value = 0.1607838857200097

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3747116666666667
$ws.Range("H4").Value = 1.124135
$ws.Range("I4").Value = 0.3914669751594584
$ws.Range("J4").Value = 0.3914669751594584
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.710753333333333
$ws.Range("N4").Value = 11.13226
$ws.Range("O4").Value = 0.2244644393883547
$ws.Range("P4").Value = 0.2244644393883547
$ws.Range("Q4").Value = 1.390462566122222
$ws.Range("R4").Value = 12.5141630951
$ws.Range("S4").Value = 0.08787041511822281
$ws.Range("T4").Value = 0.08787041511822279

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.401547
$ws.Range("H5").Value = 1.204641
$ws.Range("I5").Value = 0.4195022558883632
$ws.Range("J5").Value = 0.4195022558883631
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.030956000000001
$ws.Range("N5").Value = 18.092868
$ws.Range("O5").Value = 0.364814105361131
$ws.Range("P5").Value = 0.3648141053611309
$ws.Range("Q5").Value = 2.421712288932
$ws.Range("R5").Value = 21.79541060038801
$ws.Range("S5").Value = 0.1530403401788894
$ws.Range("T5").Value = 0.1530403401788894

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.401547
$ws.Range("H6").Value = 1.204641
$ws.Range("I6").Value = 0.4195022558883632
$ws.Range("J6").Value = 0.4195022558883631
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.789877333333333
$ws.Range("N6").Value = 20.369632
$ws.Range("O6").Value = 0.4107214552505144
$ws.Range("P6").Value = 0.4107214552505143
$ws.Range("Q6").Value = 2.726454873568
$ws.Range("R6").Value = 24.538093862112
$ws.Range("S6").Value = 0.1722985770193422
$ws.Range("T6").Value = 0.1722985770193421

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.401547
$ws.Range("H7").Value = 1.204641
$ws.Range("I7").Value = 0.4195022558883632
$ws.Range("J7").Value = 0.4195022558883631
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.710753333333333
$ws.Range("N7").Value = 11.13226
$ws.Range("O7").Value = 0.2244644393883547
$ws.Range("P7").Value = 0.2244644393883547
$ws.Range("Q7").Value = 1.49004186874
$ws.Range("R7").Value = 13.41037681866
$ws.Range("S7").Value = 0.09416333869013156
$ws.Range("T7").Value = 0.09416333869013153

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.18094
$ws.Range("H8").Value = 0.54282
$ws.Range("I8").Value = 0.1890307689521785
$ws.Range("J8").Value = 0.1890307689521785
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.030956000000001
$ws.Range("N8").Value = 18.092868
$ws.Range("O8").Value = 0.364814105361131
$ws.Range("P8").Value = 0.3648141053611309
$ws.Range("Q8").Value = 1.09124117864
$ws.Range("R8").Value = 9.821170607760001
$ws.Range("S8").Value = 0.06896109086101565
$ws.Range("T8").Value = 0.06896109086101564

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.18094
$ws.Range("H9").Value = 0.54282
$ws.Range("I9").Value = 0.1890307689521785
$ws.Range("J9").Value = 0.1890307689521785
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.789877333333333
$ws.Range("N9").Value = 20.369632
$ws.Range("O9").Value = 0.4107214552505144
$ws.Range("P9").Value = 0.4107214552505143
$ws.Range("Q9").Value = 1.228560404693333
$ws.Range("R9").Value = 11.05704364224
$ws.Range("S9").Value = 0.07763899251116251
$ws.Range("T9").Value = 0.07763899251116248

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.18094
$ws.Range("H10").Value = 0.54282
$ws.Range("I10").Value = 0.1890307689521785
$ws.Range("J10").Value = 0.1890307689521785
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.710753333333333
$ws.Range("N10").Value = 11.13226
$ws.Range("O10").Value = 0.2244644393883547
$ws.Range("P10").Value = 0.2244644393883547
$ws.Range("Q10").Value = 0.6714237081333332
$ws.Range("R10").Value = 6.042813373199999
$ws.Range("S10").Value = 0.04243068558000036
$ws.Range("T10").Value = 0.04243068558000034
